$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 is "Task 9: Register Backend" / "In Development".
# Add the new note and the date it was last touched, matching how the
# existing text-date cells (D2/D3, "01/30/2024") are stored: plain text,
# not an auto-converted date serial. Writing the date directly would get
# auto-coerced into a date value/style by the recalculated cell, so we
# write it as a quoted text formula and then convert the cell to its
# static value, which leaves behind a plain shared-string cell.
$ws.Range("C10").Value = "Having a problem with reflecting auth frontend"
$ws.Range("D10").Formula = "=""02/13/2024"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)  # xlPasteValues

# The user had since clicked into C12 (next open Notes cell).
$ws.Range("C12").Select()
